$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.220.81'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.403.84'
$ws.Range("E3").Value = '  -0.69%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.44'
$ws.Range("E5").Value = '  +1.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.34'
$ws.Range("E6").Value = '  -1.82%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.588'
$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.64'
$ws.Range("E10").Value = '  -0.82%  '

$ws.Range("E11").Value = '  +1.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.346'
$ws.Range("E12").Value = '  -2.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.67'
$ws.Range("E13").Value = '  -3.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.834.46'
$ws.Range("E14").Value = '  -0.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.136.74'
$ws.Range("E15").Value = '  +0.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000138'
$ws.Range("E16").Value = '  +0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.395.23'
$ws.Range("E17").Value = '  -0.64%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.16'
$ws.Range("E18").Value = '  -2.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.53'
$ws.Range("E19").Value = '  +2.83%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.76'
$ws.Range("E20").Value = '  -1.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.79'
$ws.Range("E21").Value = '  +1.52%  '

$ws.Range("E22").Value = '  -0.08%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '64.41'
$ws.Range("E23").Value = '  -3.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.174'
$ws.Range("E24").Value = '  +0.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.52'
$ws.Range("E25").Value = '  -2.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.11%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.38'
$ws.Range("E27").Value = '  +0.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.81'
$ws.Range("E28").Value = '  +1.75%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0769'
$ws.Range("E29").Value = '  -1.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.80'
$ws.Range("E30").Value = '  +0.93%  '

$ws.Range("E31").Value = '  +0.13%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.09'
$ws.Range("E32").Value = '  +6.94%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.402'
$ws.Range("E33").Value = '  -1.88%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.30'
$ws.Range("E34").Value = '  -2.29%  '

$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.33'
$ws.Range("E36").Value = '  +1.89%  '

$ws.Range("E37").Value = '  +0.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.18'
$ws.Range("E38").Value = '  -0.95%  '

$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '324.05'
$ws.Range("E39").Value = '  +3.12%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.60'
$ws.Range("E40").Value = '  -0.75%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.77'
$ws.Range("E41").Value = '  -2.08%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.57'
$ws.Range("E42").Value = '  +6.32%  '

$ws.Range("E43").Value = '  -3.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0971'
$ws.Range("E44").Value = '  +0.37%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.90'
$ws.Range("E45").Value = '  +1.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0515'
$ws.Range("E46").Value = '  -1.00%  '

$ws.Range("E47").Value = '  -0.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0221'
$ws.Range("E48").Value = '  -1.76%  '

$ws.Range("E49").Value = '  +0.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.57'
$ws.Range("E50").Value = '  -1.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.68'
$ws.Range("E51").Value = '  -0.08%  '
